$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "moy_pro" column header (D1), copying the header style from C1 ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "moy_pro"

# --- Copy the data-cell style (from C2) down the whole new D column ---
$ws.Range("C2").Copy()
$ws.Range("D2:D16").PasteSpecial(-4122)

# --- Update existing B/C values and populate the new D column ---
$ws.Range("B2").Value = 41.6
$ws.Range("C2").Value = 37.6
$ws.Range("D2").Value = 32.4

$ws.Range("B3").Value = 21.3
$ws.Range("C3").Value = 30.4
$ws.Range("D3").Value = 35.7

$ws.Range("B4").Value = 4.8
$ws.Range("C4").Value = 5.9
$ws.Range("D4").Value = 8.3

$ws.Range("B5").Value = 6.8
$ws.Range("C5").Value = 8.4
$ws.Range("D5").Value = 11.2

$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 20.4
$ws.Range("D6").Value = 26.7

$ws.Range("B7").Value = 0.1
$ws.Range("C7").Value = 0.1
$ws.Range("D7").Value = 0.1

$ws.Range("B8").Value = 4.8
$ws.Range("C8").Value = 5.2
$ws.Range("D8").Value = 4.6

$ws.Range("B9").Value = 14.6
$ws.Range("C9").Value = 17.2
$ws.Range("D9").Value = 27.2

$ws.Range("B10").Value = 2.8
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 2

$ws.Range("B11").Value = 8.3
$ws.Range("C11").Value = 8.7
$ws.Range("D11").Value = 20.4

$ws.Range("B12").Value = 12.1
$ws.Range("C12").Value = 11.9
$ws.Range("D12").Value = 9

# Row 13: B and C become blank (keep style), D gets a value
$ws.Range("B13").Value = 1
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 1
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = 20.2

$ws.Range("B14").Value = 10
$ws.Range("C14").Value = 20
# Row 14: D becomes blank (keep style, no value)
$ws.Range("D14").Value = 1
$ws.Range("D14").ClearContents()

$ws.Range("B15").Value = 10.7
$ws.Range("C15").Value = 7.1
$ws.Range("D15").Value = 5.1

$ws.Range("B16").Value = 5.1
$ws.Range("C16").Value = 4.5
$ws.Range("D16").Value = 3

# --- Remove the custom widths on columns B and C so they fall back to default ---
$ws.Columns("B:C").ColumnWidth = $ws.Columns("A").ColumnWidth

# --- Update the selected cell shown when the workbook is reopened ---
$ws.Range("F11").Select()
